{"js": "// Update the worksheet date and all two-digit multiplication problems\n// to the values introduced by the target commit. Each old value is\n// unique within the document, so a straightforward search-and-replace\n// per pair is sufficient and unambiguous.\nconst replacements = [\n  [\"2023-09-27 Wednesday\", \"2023-09-28 Thursday\"],\n  [\"70\u00d723=\", \"22\u00d739=\"],\n  [\"89\u00d755=\", \"28\u00d757=\"],\n  [\"37\u00d799=\", \"89\u00d786=\"],\n  [\"73\u00d721=\", \"81\u00d763=\"],\n  [\"56\u00d733=\", \"15\u00d789=\"],\n  [\"21\u00d758=\", \"34\u00d765=\"],\n  [\"20\u00d789=\", \"84\u00d755=\"],\n  [\"83\u00d721=\", \"51\u00d715=\"],\n  [\"61\u00d737=\", \"21\u00d789=\"],\n  [\"28\u00d751=\", \"30\u00d744=\"],\n  [\"82\u00d736=\", \"60\u00d733=\"],\n  [\"31\u00d780=\", \"48\u00d754=\"],\n  [\"39\u00d785=\", \"31\u00d799=\"],\n  [\"74\u00d790=\", \"49\u00d761=\"],\n  [\"68\u00d784=\", \"28\u00d753=\"],\n  [\"82\u00d778=\", \"85\u00d758=\"],\n  [\"89\u00d734=\", \"72\u00d794=\"],\n  [\"75\u00d781=\", \"76\u00d782=\"],\n  [\"20\u00d738=\", \"62\u00d743=\"],\n  [\"64\u00d743=\", \"90\u00d791=\"],\n  [\"64\u00d718=\", \"97\u00d766=\"],\n  [\"12\u00d733=\", \"60\u00d726=\"],\n  [\"84\u00d769=\", \"95\u00d748=\"],\n  [\"12\u00d770=\", \"84\u00d778=\"],\n  [\"86\u00d795=\", \"59\u00d774=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2023-09-27 Wednesday\"; New = \"2023-09-28 Thursday\" },\n    @{ Old = \"70\u00d723=\"; New = \"22\u00d739=\" },\n    @{ Old = \"89\u00d755=\"; New = \"28\u00d757=\" },\n    @{ Old = \"37\u00d799=\"; New = \"89\u00d786=\" },\n    @{ Old = \"73\u00d721=\"; New = \"81\u00d763=\" },\n    @{ Old = \"56\u00d733=\"; New = \"15\u00d789=\" },\n    @{ Old = \"21\u00d758=\"; New = \"34\u00d765=\" },\n    @{ Old = \"20\u00d789=\"; New = \"84\u00d755=\" },\n    @{ Old = \"83\u00d721=\"; New = \"51\u00d715=\" },\n    @{ Old = \"61\u00d737=\"; New = \"21\u00d789=\" },\n    @{ Old = \"28\u00d751=\"; New = \"30\u00d744=\" },\n    @{ Old = \"82\u00d736=\"; New = \"60\u00d733=\" },\n    @{ Old = \"31\u00d780=\"; New = \"48\u00d754=\" },\n    @{ Old = \"39\u00d785=\"; New = \"31\u00d799=\" },\n    @{ Old = \"74\u00d790=\"; New = \"49\u00d761=\" },\n    @{ Old = \"68\u00d784=\"; New = \"28\u00d753=\" },\n    @{ Old = \"82\u00d778=\"; New = \"85\u00d758=\" },\n    @{ Old = \"89\u00d734=\"; New = \"72\u00d794=\" },\n    @{ Old = \"75\u00d781=\"; New = \"76\u00d782=\" },\n    @{ Old = \"20\u00d738=\"; New = \"62\u00d743=\" },\n    @{ Old = \"64\u00d743=\"; New = \"90\u00d791=\" },\n    @{ Old = \"64\u00d718=\"; New = \"97\u00d766=\" },\n    @{ Old = \"12\u00d733=\"; New = \"60\u00d726=\" },\n    @{ Old = \"84\u00d769=\"; New = \"95\u00d748=\" },\n    @{ Old = \"12\u00d770=\"; New = \"84\u00d778=\" },\n    @{ Old = \"86\u00d795=\"; New = \"59\u00d774=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"Could not find text to replace: $($pair.Old)\"\n    }\n}\n"}
